$d = $word.ActiveDocument

$r1 = $d.Content
$r1.Find.Execute(" is 2.4 ", $false, $false, $false, $false, $false, $true, 1, $false, " is 2.45 ", 2)

$r2 = $d.Content
$r2.Find.Execute(" 2.4 ", $false, $false, $false, $false, $false, $true, 1, $false, " 2.45 ", 2)
